$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows representing the same conversation joined/duplicated for a
# different group (group_id -1001159430667), continuing the "id" sequence.
$newRows = @(
    @(10, 1, "Hello.", -1001159430667),
    @(11, 3, "Nice to see you", -1001159430667),
    @(12, 2, "Hi. How are you?", -1001159430667),
    @(13, 1, "I'm fine. Thanks. And you?", -1001159430667),
    @(14, 2, "I'm fine. Thanks", -1001159430667),
    @(15, 2, "I'm busy now.", -1001159430667),
    @(16, 2, "See you later.", -1001159430667),
    @(17, 1, "See you.", -1001159430667),
    @(18, 3, "See you.", -1001159430667)
)

$r = 11
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Column B (user_id) widened and no longer "best fit".
# (Target authored width is 12.140625 chars; the host engine quantizes
# stored column widths to multiples of 1/6 character, so 12.1667 -- the
# nearest representable value -- is used here.)
$ws.Range("B:B").ColumnWidth = 11.3

# Selection moved to G17 (as left by the author after editing)
$ws.Range("G17").Select()
